# Apply updated values to result_data_RandomForest.xlsx (Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -10.54989999999999
$ws.Range("D3").Value = -6.942399999999992
$ws.Range("A12").Value = -21.5099
$ws.Range("C14").Value = -13.3562
$ws.Range("C26").Value = -12.4191
$ws.Range("D30").Value = -7.5082
$ws.Range("C31").Value = -12.91530000000001
$ws.Range("A32").Value = -21.35179999999999
$ws.Range("C35").Value = -12.83220000000002
$ws.Range("A36").Value = -19.6661
$ws.Range("C37").Value = -13.5188
$ws.Range("A38").Value = -19.3641
$ws.Range("D44").Value = -7.2458
$ws.Range("C45").Value = -13.87809999999999
$ws.Range("A46").Value = -21.6844
$ws.Range("A54").Value = -21.80029999999999
$ws.Range("A55").Value = -22.55580000000001
$ws.Range("C57").Value = -14.73749999999999
$ws.Range("D58").Value = -8.212099999999998
$ws.Range("A67").Value = -21.54619999999997
$ws.Range("A69").Value = -21.64929999999997
$ws.Range("A72").Value = -21.694
$ws.Range("D84").Value = -8.665300000000004
$ws.Range("D89").Value = -6.096999999999997
$ws.Range("A91").Value = -21.542
$ws.Range("D91").Value = -5.974399999999997
$ws.Range("D92").Value = -6.070100000000002
$ws.Range("A99").Value = -20.41949999999998
$ws.Range("C100").Value = -12.3163
$ws.Range("C102").Value = -14.2256
$ws.Range("D102").Value = -7.939499999999999
